$wb = $excel.ActiveWorkbook

# --- Sheet "Data Set 0 Timings (Pd)" ---
$wsPd = $wb.Worksheets.Item("Data Set 0 Timings (Pd)")

$wsPd.Cells.Item(2, 6).Value = 15.106
$wsPd.Cells.Item(2, 7).Value = 0.284
$wsPd.Cells.Item(2, 8).Value = 9.864000000000001
$wsPd.Cells.Item(2, 9).Value = 2.275
$wsPd.Cells.Item(2, 10).Value = 2.32
$wsPd.Cells.Item(2, 12).Value = 0.00021

$wsPd.Cells.Item(3, 6).Value = 15.781
$wsPd.Cells.Item(3, 7).Value = 0.308
$wsPd.Cells.Item(3, 8).Value = 10.459
$wsPd.Cells.Item(3, 9).Value = 2.349
$wsPd.Cells.Item(3, 10).Value = 2.381
$wsPd.Cells.Item(3, 12).Value = 0.00021

$wsPd.Cells.Item(4, 6).Value = 18.494
$wsPd.Cells.Item(4, 7).Value = 0.345
$wsPd.Cells.Item(4, 8).Value = 12.481
$wsPd.Cells.Item(4, 9).Value = 2.679
$wsPd.Cells.Item(4, 10).Value = 2.689
$wsPd.Cells.Item(4, 12).Value = 0.00024

$wsPd.Cells.Item(5, 6).Value = 17.479
$wsPd.Cells.Item(5, 7).Value = 0.35
$wsPd.Cells.Item(5, 8).Value = 11.618
$wsPd.Cells.Item(5, 9).Value = 2.566
$wsPd.Cells.Item(5, 10).Value = 2.593
$wsPd.Cells.Item(5, 12).Value = 0.00023

$wsPd.Cells.Item(6, 6).Value = 17.426
$wsPd.Cells.Item(6, 7).Value = 0.337
$wsPd.Cells.Item(6, 8).Value = 11.511
$wsPd.Cells.Item(6, 9).Value = 2.605
$wsPd.Cells.Item(6, 10).Value = 2.619
$wsPd.Cells.Item(6, 12).Value = 0.00024

# --- Sheet "Data Set 0 Timings (TD)" ---
$wsTd = $wb.Worksheets.Item("Data Set 0 Timings (TD)")

$wsTd.Cells.Item(2, 6).Value = 30.998
$wsTd.Cells.Item(2, 7).Value = 0.257
$wsTd.Cells.Item(2, 8).Value = 0.338
$wsTd.Cells.Item(2, 9).Value = 11.55
$wsTd.Cells.Item(2, 10).Value = 11.428
$wsTd.Cells.Item(2, 12).Value = 0.00023

$wsTd.Cells.Item(3, 6).Value = 30.864
$wsTd.Cells.Item(3, 7).Value = 0.258
$wsTd.Cells.Item(3, 8).Value = 0.339
$wsTd.Cells.Item(3, 9).Value = 11.644
$wsTd.Cells.Item(3, 10).Value = 11.653
$wsTd.Cells.Item(3, 12).Value = 0.00024

$wsTd.Cells.Item(4, 6).Value = 32.022
$wsTd.Cells.Item(4, 7).Value = 0.254
$wsTd.Cells.Item(4, 8).Value = 0.338
$wsTd.Cells.Item(4, 9).Value = 12.816
$wsTd.Cells.Item(4, 10).Value = 11.519
$wsTd.Cells.Item(4, 12).Value = 0.00023

$wsTd.Cells.Item(5, 6).Value = 31.021
$wsTd.Cells.Item(5, 7).Value = 0.256
$wsTd.Cells.Item(5, 8).Value = 0.339
$wsTd.Cells.Item(5, 9).Value = 11.591
$wsTd.Cells.Item(5, 10).Value = 11.66
$wsTd.Cells.Item(5, 12).Value = 0.00024

$wsTd.Cells.Item(6, 6).Value = 31.013
$wsTd.Cells.Item(6, 7).Value = 0.256
$wsTd.Cells.Item(6, 8).Value = 0.338
$wsTd.Cells.Item(6, 9).Value = 11.727
$wsTd.Cells.Item(6, 10).Value = 11.576
$wsTd.Cells.Item(6, 12).Value = 0.00024

# --- Sheet "Data Set 0 Timings (combined)" ---
# Rows 2-6 mirror the "Pd" sheet values, rows 7-11 mirror the "TD" sheet values.
$wsComb = $wb.Worksheets.Item("Data Set 0 Timings (combined)")

$wsComb.Cells.Item(2, 6).Value = 15.106
$wsComb.Cells.Item(2, 7).Value = 0.284
$wsComb.Cells.Item(2, 8).Value = 9.864000000000001
$wsComb.Cells.Item(2, 9).Value = 2.275
$wsComb.Cells.Item(2, 10).Value = 2.32
$wsComb.Cells.Item(2, 12).Value = 0.00021

$wsComb.Cells.Item(3, 6).Value = 15.781
$wsComb.Cells.Item(3, 7).Value = 0.308
$wsComb.Cells.Item(3, 8).Value = 10.459
$wsComb.Cells.Item(3, 9).Value = 2.349
$wsComb.Cells.Item(3, 10).Value = 2.381
$wsComb.Cells.Item(3, 12).Value = 0.00021

$wsComb.Cells.Item(4, 6).Value = 18.494
$wsComb.Cells.Item(4, 7).Value = 0.345
$wsComb.Cells.Item(4, 8).Value = 12.481
$wsComb.Cells.Item(4, 9).Value = 2.679
$wsComb.Cells.Item(4, 10).Value = 2.689
$wsComb.Cells.Item(4, 12).Value = 0.00024

$wsComb.Cells.Item(5, 6).Value = 17.479
$wsComb.Cells.Item(5, 7).Value = 0.35
$wsComb.Cells.Item(5, 8).Value = 11.618
$wsComb.Cells.Item(5, 9).Value = 2.566
$wsComb.Cells.Item(5, 10).Value = 2.593
$wsComb.Cells.Item(5, 12).Value = 0.00023

$wsComb.Cells.Item(6, 6).Value = 17.426
$wsComb.Cells.Item(6, 7).Value = 0.337
$wsComb.Cells.Item(6, 8).Value = 11.511
$wsComb.Cells.Item(6, 9).Value = 2.605
$wsComb.Cells.Item(6, 10).Value = 2.619
$wsComb.Cells.Item(6, 12).Value = 0.00024

$wsComb.Cells.Item(7, 6).Value = 30.998
$wsComb.Cells.Item(7, 7).Value = 0.257
$wsComb.Cells.Item(7, 8).Value = 0.338
$wsComb.Cells.Item(7, 9).Value = 11.55
$wsComb.Cells.Item(7, 10).Value = 11.428
$wsComb.Cells.Item(7, 12).Value = 0.00023

$wsComb.Cells.Item(8, 6).Value = 30.864
$wsComb.Cells.Item(8, 7).Value = 0.258
$wsComb.Cells.Item(8, 8).Value = 0.339
$wsComb.Cells.Item(8, 9).Value = 11.644
$wsComb.Cells.Item(8, 10).Value = 11.653
$wsComb.Cells.Item(8, 12).Value = 0.00024

$wsComb.Cells.Item(9, 6).Value = 32.022
$wsComb.Cells.Item(9, 7).Value = 0.254
$wsComb.Cells.Item(9, 8).Value = 0.338
$wsComb.Cells.Item(9, 9).Value = 12.816
$wsComb.Cells.Item(9, 10).Value = 11.519
$wsComb.Cells.Item(9, 12).Value = 0.00023

$wsComb.Cells.Item(10, 6).Value = 31.021
$wsComb.Cells.Item(10, 7).Value = 0.256
$wsComb.Cells.Item(10, 8).Value = 0.339
$wsComb.Cells.Item(10, 9).Value = 11.591
$wsComb.Cells.Item(10, 10).Value = 11.66
$wsComb.Cells.Item(10, 12).Value = 0.00024

$wsComb.Cells.Item(11, 6).Value = 31.013
$wsComb.Cells.Item(11, 7).Value = 0.256
$wsComb.Cells.Item(11, 8).Value = 0.338
$wsComb.Cells.Item(11, 9).Value = 11.727
$wsComb.Cells.Item(11, 10).Value = 11.576
$wsComb.Cells.Item(11, 12).Value = 0.00024
